$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new inventory row 75 (no description / column C left blank)
$ws.Range("A75").Value = "H1R39U"
$ws.Range("B75").Value = "Goma separador de bandeja Samsung"
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 10000
$ws.Range("F75").Value = 72
$ws.Range("G75").Value = 0
$ws.Range("H75").Formula = "=(E75-D75)*G75"
$ws.Range("I75").Formula = "=D75*F75"
$ws.Range("J75").Value = 0
